$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# Copy the Display value (column C) into the empty Definition column (column D)
# for each concept row.
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($row, 3).Value2
}
